$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1487
$ws1.Range("F9").Value = 276

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1487
$ws4.Range("F9").Value = 276
